$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.718.28'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '1.889.39'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4731'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2925'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06530'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07801'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").Value = '1.888.45'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7374'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.253'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.66%  '
$ws.Range("D17").Value = '30.705.71'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007542'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '2.136.43'
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.320'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.251'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.227'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.918'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.339'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.94%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09742'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.491'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.300'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.193'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.126'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6968'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.727'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("E39").Value = '  +2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.340'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '76.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.000'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4277'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8367'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.536'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.49%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.046'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '914.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05756'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.02%  '
